$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 72
$ws.Range("I2").Value = 169
$ws.Range("J2").Value = 702
$ws.Range("K2").Value = 5
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 120
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 10
$ws.Range("T2").Value = 120
$ws.Range("V2").Value = 1083
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1112
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 6
